$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Installation of jsonwebtoken : npm install jsonwebtoken" paragraph
#    -> split into runs with spell-check w:proofErr markers bracketing
#       the non-dictionary words "npm" and "jsonwebtoken".
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(20)
$r1 = $p1.Range
$range1 = $d.Range($r1.Start, $r1.End - 1)

$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Installation of jsonwebtoken : </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> install </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jsonwebtoken</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$range1.InsertXML($xml1)

# ------------------------------------------------------------------
# 2) "Make a file named user.js ... decrypted back." paragraph
#    -> split into runs with w:proofErr markers bracketing
#       "materials." + "hashing" (one joint span, since they were
#       adjacent with no space) and "datas".
#       <w:lastRenderedPageBreak/> stays on the first run.
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(23)
$r2 = $p2.Range
$range2 = $d.Range($r2.Start, $r2.End - 1)

$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Make a file named user.js and copy the code from the study </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>materials.</w:t></w:r><w:r><w:t>hashing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is the one way where data can be encrypted but it’s a one way . </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>datas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> that has been encrypted cannot be decrypted back.</w:t></w:r></w:p>
'@
$range2.InsertXML($xml2)

# ------------------------------------------------------------------
# 3) New paragraph appended right after it: "In the blog.js set the
#    reference of the user for setting the scema ." with a w:proofErr
#    span around "scema".
#    A clean paragraph break is created first (InsertAfter with a
#    literal carriage return), then the new paragraph's runs are
#    replaced via InsertXML so no stray empty run is left behind.
# ------------------------------------------------------------------
$p2b = $d.Paragraphs.Item(23)
$r2b = $p2b.Range
$breakPoint = $d.Range($r2b.End - 1, $r2b.End - 1)
$breakPoint.InsertAfter("`rPLACEHOLDER_NEW_PARA")

$p3 = $d.Paragraphs.Item(24)
$r3 = $p3.Range
$range3 = $d.Range($r3.Start, $r3.End - 1)

$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">In the blog.js set the reference of the user for setting the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>scema</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> .</w:t></w:r></w:p>
'@
$range3.InsertXML($xml3)
